$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Simple text replacements (value only, style unchanged) ---
$ws.Range("C1").Value  = "Operator Controller"
$ws.Range("H1").Value  = "v0.5 (Mar 2, 2020)"
$ws.Range("B3").Value  = " (hold) Vision Align and Shoot Close"
$ws.Range("H3").Value  = "Vision Align and Shoot Far (hold)"
$ws.Range("B4").Value  = "(hold )Shoot Close"
$ws.Range("H4").Value  = "Shoot Far (hold)"
$ws.Range("H6").Value  = "Control Panel Arm Fold/Unfold"
$ws.Range("H7").Value  = "Expel (hold)"
$ws.Range("H8").Value  = "Intake (hold)"
$ws.Range("H9").Value  = "*[see Shift table]"
$ws.Range("E16").Value = "Climb Arm Lock"
$ws.Range("C21").Value = "Driver Controller"
$ws.Range("H23").Value = "Disable Auto-Shift (hold)"
$ws.Range("B24").Value = "(hold) Climb Adjust Left"
$ws.Range("H24").Value = "Climb Adjust Right (hold)"

# --- B10 gains the text that used to live in B11; style goes from the
#     "blank placeholder" (s=4, grey) to the "filled" look (s=3) that the
#     other populated rows in this column use ---
$ws.Range("B10").Value = "(Control Panel Rotate)"
$ws.Range("B10").Font.Color = $ws.Range("B9").Font.Color
$ws.Range("B10").Font.ColorIndex = $ws.Range("B9").Font.ColorIndex

# --- B11 is fully cleared out: no text, no special formatting left behind ---
$ws.Range("B11").Value = ""
$ws.Range("B11").Style = "Normal"

# --- H16 text is cleared, and the leftover yellow highlight fill is
#     removed, while keeping the bold/centered look of the cell ---
$ws.Range("H16").Value = ""
$ws.Range("H16").Style = "Normal"
$ws.Range("H16").Font.Bold = $true
$ws.Range("H16").Font.Size = 10
$ws.Range("H16").HorizontalAlignment = -4108

# --- Selection moves from H18 to B25 ---
$ws.Range("B25").Select()
